$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 215, shifting existing rows 215:300 down to 216:301
$ws.Rows.Item(215).Insert()

# Populate the newly inserted row 215 with the new record's data
$ws.Range("A215").Value = 11
$ws.Range("B215").Value = "Vega Monumental Concepción"
$ws.Range("C215").Value = "Bíobío"
$ws.Range("D215").Value = 44510
$ws.Range("E215").Value = 8
$ws.Range("F215").Value = "Fruta"
$ws.Range("G215").Value = 100102
$ws.Range("H215").Value = "Cítricos"
$ws.Range("I215").Value = 100102003
$ws.Range("J215").Value = "Limón"
$ws.Range("K215").Value = "Sin especificar"
$ws.Range("L215").Value = "1a amarillo"
$ws.Range("M215").Value = 350
$ws.Range("N215").Value = 6000
$ws.Range("O215").Value = 7000
$ws.Range("P215").Value = 6571
$ws.Range("Q215").Value = "`$/malla 16 kilos"
$ws.Range("R215").Value = "Región de O'Higgins"
$ws.Range("S215").Value = 411
$ws.Range("T215").Value = 16
